# ------------------------------------------------------------------
# feat: add 2022-Q1 data
#  - inserts a new worksheet "2022-Q1" (with fund holdings detail)
#    right before the "总计" summary sheet
#  - adds a new top row to "总计" summarizing the 2022-Q1 quarter
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" worksheet, placed right before "总计" ---
$refSheet = $wb.Worksheets.Item("2021-Q4")
$totalRef = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalRef)
$ws.Name = "2022-Q1"

# NOTE: inserting a sheet shifts tab positions, and sheet variables here
# resolve by position - re-fetch "总计" by name now that it has moved so
# $total reliably points at the summary sheet (not the new one).
$total = $wb.Worksheets.Item("总计")

# Copy header-row (B1:H1) and index-column (A2:A11) formatting from the
# "2021-Q4" sheet, which already has the right style (bold header, bordered
# index column) for this table layout.
$refSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2:A11").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Index column (A2:A11) = 0..9
for ($i = 2; $i -le 11; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# Columns B-G hold text values (fund codes/names/numbers-as-text, matching
# the source data which stores these as text, not numbers, e.g. to keep
# leading zeros in fund codes). Force text storage via NumberFormat, then
# ClearFormats() afterwards to drop back to the default (unstyled) cell
# format without disturbing the text values already entered.
$textRange = $ws.Range("B2:G11")
$textRange.NumberFormat = "@"

$ws.Range("B2").Value = "007449"
$ws.Range("C2").Value = "兴全多维价值混合A"
$ws.Range("D2").Value = "28.03"
$ws.Range("E2").Value = "84.40"
$ws.Range("F2").Value = "3.24"
$ws.Range("G2").Value = "0.9082"
$ws.Range("B3").Value = "010751"
$ws.Range("C3").Value = "宝盈优质成长混合A"
$ws.Range("D3").Value = "5.64"
$ws.Range("E3").Value = "92.80"
$ws.Range("F3").Value = "4.37"
$ws.Range("G3").Value = "0.2465"
$ws.Range("B4").Value = "001543"
$ws.Range("C4").Value = "宝盈新锐灵活配置混合A"
$ws.Range("D4").Value = "3.21"
$ws.Range("E4").Value = "93.26"
$ws.Range("F4").Value = "4.68"
$ws.Range("G4").Value = "0.1502"
$ws.Range("B5").Value = "007450"
$ws.Range("C5").Value = "兴全多维价值混合C"
$ws.Range("D5").Value = "4.48"
$ws.Range("E5").Value = "84.40"
$ws.Range("F5").Value = "3.24"
$ws.Range("G5").Value = "0.1452"
$ws.Range("B6").Value = "011550"
$ws.Range("C6").Value = "湘财创新成长一年持有期混合A"
$ws.Range("D6").Value = "2.62"
$ws.Range("E6").Value = "93.51"
$ws.Range("F6").Value = "4.44"
$ws.Range("G6").Value = "0.1163"
$ws.Range("B7").Value = "010752"
$ws.Range("C7").Value = "宝盈优质成长混合C"
$ws.Range("D7").Value = "0.78"
$ws.Range("E7").Value = "92.80"
$ws.Range("F7").Value = "4.37"
$ws.Range("G7").Value = "0.0341"
$ws.Range("B8").Value = "010076"
$ws.Range("C8").Value = "湘财长弘灵活配置混合A"
$ws.Range("D8").Value = "0.32"
$ws.Range("E8").Value = "91.72"
$ws.Range("F8").Value = "5.20"
$ws.Range("G8").Value = "0.0166"
$ws.Range("B9").Value = "011551"
$ws.Range("C9").Value = "湘财创新成长一年持有期混合C"
$ws.Range("D9").Value = "0.28"
$ws.Range("E9").Value = "93.51"
$ws.Range("F9").Value = "4.44"
$ws.Range("G9").Value = "0.0124"
$ws.Range("B10").Value = "007578"
$ws.Range("C10").Value = "宝盈新锐灵活配置混合C"
$ws.Range("D10").Value = "0.20"
$ws.Range("E10").Value = "93.26"
$ws.Range("F10").Value = "4.68"
$ws.Range("G10").Value = "0.0094"
$ws.Range("B11").Value = "010077"
$ws.Range("C11").Value = "湘财长弘灵活配置混合C"
$ws.Range("D11").Value = "0.10"
$ws.Range("E11").Value = "91.72"
$ws.Range("F11").Value = "5.20"
$ws.Range("G11").Value = "0.0052"

# Drop back to the default (unstyled) format now that the text values are
# set, so these cells end up with no explicit style - matching the rest of
# the workbook's data rows.
$textRange.ClearFormats()

# H column (仓位排名) holds real numbers.
$ws.Range("H2").Value = 4
$ws.Range("H3").Value = 7
$ws.Range("H4").Value = 8
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 5
$ws.Range("H7").Value = 7
$ws.Range("H8").Value = 4
$ws.Range("H9").Value = 5
$ws.Range("H10").Value = 8
$ws.Range("H11").Value = 4

# --- 2. Insert a new top data-row in "总计" for the 2022-Q1 quarter ---
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

# Re-apply the index-column style (A2) from the row below it, matching the
# rest of the index column.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 1.64

# Renumber the index column for the rows pushed down by the insert.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
